$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.640.16"
$ws.Range("E2").Value = "  -1.18%  "

# Row 3
$ws.Range("D3").Value = "3.862.65"
$ws.Range("E3").Value = "  -2.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.73%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.712"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.36%  "

# Row 10
$ws.Range("E10").Value = "  -6.18%  "

# Row 11
$ws.Range("E11").Value = "  -8.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.57"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.35"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.23%  "

# Row 14
$ws.Range("D14").Value = "4.484.63"
$ws.Range("E14").Value = "  -1.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.44%  "

# Row 16
$ws.Range("D16").Value = "3.859.40"
$ws.Range("E16").Value = "  -2.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.76%  "

# Row 18
$ws.Range("E18").Value = "  -2.18%  "

# Row 19
$ws.Range("E19").Value = "  +2.53%  "

# Row 20
$ws.Range("D20").Value = "68.670.96"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.99%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.94%  "

# Row 24
$ws.Range("E24").Value = "  -3.01%  "

# Row 25
$ws.Range("E25").Value = "  +5.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.38%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.50"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.74%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.25"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "679.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.28%  "

# Row 31
$ws.Range("E31").Value = "  -5.36%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.89"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +14.50%  "

# Row 33
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.75%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.453"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.48%  "

# Row 36
$ws.Range("B36").Value = "ThetaToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.54"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +15.59%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.60"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.38%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0829"
$ws.Range("E38").Value = "  -8.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.148"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.07%  "

# Row 40
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0474"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.18%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.13"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.35%  "

# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.85%  "

# Row 46
$ws.Range("E46").Value = "  -1.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.76%  "

# Row 48
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000269"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +11.91%  "

# Row 49
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.28"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.94%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0340"
$ws.Range("E50").Value = "  -5.59%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.92"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.32%  "
